# Append a new row (row 28) to Sheet1 with the next day's gold-price entry,
# mirroring the existing "Date" / "Gold data" column layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: new date label for the appended row.
$ws.Cells.Item(28, 1).Value = "19-10-2025"

# Column B: same gold-price text as the prior day's row (row 27), copied via
# .Text so the existing shared string is reused instead of creating a new one.
$ws.Cells.Item(28, 2).Value = $ws.Cells.Item(27, 2).Text
